$d = $word.ActiveDocument

# Second table in the document holds the per-day travel/activity log.
$t2 = $d.Tables.Item(2)

# Row 3 (14/05/2022 entry): total-hours cell "26" -> "2:30"
$hoursCell1 = $t2.Rows.Item(3).Cells.Item(13)
$hoursCell1.Range.Text = "2:30"

# Row 4 (16/05/2022 entry): total-hours cell "4" -> "4:40"
$hoursCell2 = $t2.Rows.Item(4).Cells.Item(13)
$hoursCell2.Range.Text = "4:40"

# Remove one of the trailing fully-blank rows (report finished, one fewer
# spare row is needed).
$t2.Rows.Item(5).Delete()
